$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").Value = 13
$ws.Range("B14").Formula = "=(4/10)*10"
$ws.Range("C14").Formula = "=(5/12)*10"
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = 10
$ws.Range("G14").Value = 8.5
$ws.Range("H14").Value = 5
$ws.Range("I14").Value = 5
$ws.Range("J14").Formula = "=(5/17)*10"
$ws.Range("K14").Formula = "=(14/35)*10"

# Row 34
$ws.Range("A34").Value = 13
$ws.Range("B34").Formula = "=(SUM(B14:R14)/B17) * 40"
$ws.Range("C34").Value = 5
$ws.Range("D34").Value = 0
$ws.Range("E34").Formula = "=B34+C34+D34"
$ws.Range("F34").Formula = "=E34"

# Selection state, matching the saved workbook view
$ws.Range("H18").Select()
